# Atualizacao de bases das ligas, do dia: 01-04-2024 as 22:23
# Hungary NB I: insert 2 new (earlier-played) matches before the current
# row 154, refresh the closing odds of the 5 matches that shift down, and
# append one brand-new fixture at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room: push the existing rows 154-158 down to 156-160 so the two
#    new fixtures (ids 6818349 / 6818350) can be inserted in chronological
#    order.
# ---------------------------------------------------------------------
$ws.Rows.Item(154).Resize(2).Insert()

# Re-apply the same look as the rest of the table to the freshly inserted
# rows (bold/centered/bordered id column, custom date-time format column)
# instead of whatever Excel guessed from the row above.
$ws.Cells.Item(153, 1).Copy()
$ws.Range("A154:A155").PasteSpecial(-4122)
$ws.Cells.Item(153, 5).Copy()
$ws.Range("E154:E155").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Fill in the two newly inserted rows.
# ---------------------------------------------------------------------
$newRows = @(
    @(154, 152, 6818349, 45382.4375,         "Paksi",            "Ujpest",           1, 2, "A",
      1.666, 3.5,   4.333, 1.666, 3.6,  4.2,   -0.75, 1.925, 1.925, 2.75, 1.8,   2.05,
      -1, -1, 3.2, -1, 0.925, 0.4, -0.5),
    @(155, 153, 6818350, 45382.54166666666,  "Mezokovesd Zsory", "Ferencvarosi TC",  0, 3, "A",
      7.5,   4.333, 1.333, 9.5,   4.75, 1.25,  1.5,   2,     1.85,  3,    2.025, 1.825,
      -1, -1, 0.25, -1, 0.8500000000000001, 0, 0)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value  = $r[1]   # A  id
    $ws.Cells.Item($row, 2).Value  = $r[2]   # B  id
    $ws.Cells.Item($row, 3).Value  = "Hungary NB I"  # C  Div
    $ws.Cells.Item($row, 4).Value  = "Hungary NB I"  # D  Div Original Name
    $ws.Cells.Item($row, 5).Value  = $r[3]   # E  Date
    $ws.Cells.Item($row, 6).Value  = $r[4]   # F  HomeTeam
    $ws.Cells.Item($row, 7).Value  = $r[5]   # G  AwayTeam
    $ws.Cells.Item($row, 8).Value  = $r[6]   # H  FTHG
    $ws.Cells.Item($row, 9).Value  = $r[7]   # I  FTAG
    $ws.Cells.Item($row, 10).Value = $r[8]   # J  FTR
    $ws.Cells.Item($row, 11).Value = $r[9]   # K  oddH_op
    $ws.Cells.Item($row, 12).Value = $r[10]  # L  oddD_op
    $ws.Cells.Item($row, 13).Value = $r[11]  # M  oddA_op
    $ws.Cells.Item($row, 14).Value = $r[12]  # N  oddH
    $ws.Cells.Item($row, 15).Value = $r[13]  # O  oddD
    $ws.Cells.Item($row, 16).Value = $r[14]  # P  oddA
    $ws.Cells.Item($row, 17).Value = $r[15]  # Q  Ah
    $ws.Cells.Item($row, 18).Value = $r[16]  # R  oddAHH
    $ws.Cells.Item($row, 19).Value = $r[17]  # S  oddAHA
    $ws.Cells.Item($row, 20).Value = $r[18]  # T  AhOU
    $ws.Cells.Item($row, 21).Value = $r[19]  # U  oddAHOver
    $ws.Cells.Item($row, 22).Value = $r[20]  # V  oddAHUnder
    $ws.Cells.Item($row, 23).Value = $r[21]  # W  PLH
    $ws.Cells.Item($row, 24).Value = $r[22]  # X  PLD
    $ws.Cells.Item($row, 25).Value = $r[23]  # Y  PLA
    $ws.Cells.Item($row, 26).Value = $r[24]  # Z  PL_Ahh
    $ws.Cells.Item($row, 27).Value = $r[25]  # AA PL_Aha
    $ws.Cells.Item($row, 28).Value = $r[26]  # AB PL_AhOver
    $ws.Cells.Item($row, 29).Value = $r[27]  # AC PL_AhUnder
}

# ---------------------------------------------------------------------
# 3) Refresh the "id" column (A) for the 5 rows that shifted from
#    154-158 down to 156-160 (A holds row-2), and update the handful of
#    odds columns whose closing values moved since the last scrape.
# ---------------------------------------------------------------------
for ($row = 156; $row -le 160; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 2
}

$ws.Cells.Item(156, 14).Value = 2.05    # N156 oddH
$ws.Cells.Item(156, 16).Value = 3.25    # P156 oddA
$ws.Cells.Item(156, 18).Value = 1.825   # R156 oddAHH
$ws.Cells.Item(156, 19).Value = 2.025   # S156 oddAHA
$ws.Cells.Item(156, 21).Value = 1.825   # U156 oddAHOver
$ws.Cells.Item(156, 22).Value = 2.025   # V156 oddAHUnder

$ws.Cells.Item(157, 14).Value = 1.75    # N157 oddH
$ws.Cells.Item(157, 16).Value = 4.333   # P157 oddA
$ws.Cells.Item(157, 18).Value = 2.05    # R157 oddAHH
$ws.Cells.Item(157, 19).Value = 1.8     # S157 oddAHA
$ws.Cells.Item(157, 21).Value = 1.875   # U157 oddAHOver
$ws.Cells.Item(157, 22).Value = 1.975   # V157 oddAHUnder

$ws.Cells.Item(158, 18).Value = 1.85    # R158 oddAHH
$ws.Cells.Item(158, 19).Value = 2       # S158 oddAHA

$ws.Cells.Item(159, 14).Value = 1.85    # N159 oddH
$ws.Cells.Item(159, 15).Value = 3.4     # O159 oddD
$ws.Cells.Item(159, 16).Value = 3.8     # P159 oddA
$ws.Cells.Item(159, 18).Value = 1.9     # R159 oddAHH
$ws.Cells.Item(159, 19).Value = 1.95    # S159 oddAHA

$ws.Cells.Item(160, 14).Value = 1.5     # N160 oddH
$ws.Cells.Item(160, 16).Value = 5.75    # P160 oddA
$ws.Cells.Item(160, 18).Value = 1.825   # R160 oddAHH
$ws.Cells.Item(160, 19).Value = 2.025   # S160 oddAHA
$ws.Cells.Item(160, 21).Value = 1.925   # U160 oddAHOver
$ws.Cells.Item(160, 22).Value = 1.925   # V160 oddAHUnder

# ---------------------------------------------------------------------
# 4) Append the brand-new fixture (id 6818357) as row 161, copying the
#    look of the row above it for the id/date columns.
# ---------------------------------------------------------------------
$ws.Cells.Item(160, 1).Copy()
$ws.Cells.Item(161, 1).PasteSpecial(-4122)
$ws.Cells.Item(160, 5).Copy()
$ws.Cells.Item(161, 5).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(161, 1).Value  = 159
$ws.Cells.Item(161, 2).Value  = 6818357
$ws.Cells.Item(161, 3).Value  = "Hungary NB I"
$ws.Cells.Item(161, 4).Value  = "Hungary NB I"
$ws.Cells.Item(161, 5).Value  = 45389.52083333334
$ws.Cells.Item(161, 6).Value  = "Kecskemeti TE"
$ws.Cells.Item(161, 7).Value  = "Puskas Academy"
$ws.Cells.Item(161, 11).Value = 2.4
$ws.Cells.Item(161, 12).Value = 3.2
$ws.Cells.Item(161, 13).Value = 2.6
$ws.Cells.Item(161, 14).Value = 3
$ws.Cells.Item(161, 15).Value = 3.25
$ws.Cells.Item(161, 16).Value = 2.15
$ws.Cells.Item(161, 17).Value = 0.25
$ws.Cells.Item(161, 18).Value = 1.925
$ws.Cells.Item(161, 19).Value = 1.925
$ws.Cells.Item(161, 20).Value = 2.5
$ws.Cells.Item(161, 21).Value = 1.975
$ws.Cells.Item(161, 22).Value = 1.875
$ws.Cells.Item(161, 23).Value = 0
$ws.Cells.Item(161, 24).Value = 0
$ws.Cells.Item(161, 25).Value = 0
$ws.Cells.Item(161, 26).Value = 0
$ws.Cells.Item(161, 27).Value = 0
